# Auto-generated edit script applying cosinor re-run values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = '[0.3034143856664677, 13.0229250603252]'
$ws.Range("N2").Value = 0.04043375024012086
$ws.Range("O2").Value = 0.04043375024012086
$ws.Range("U2").Value = '[5.096759098678909, 12.951271220098285]'
$ws.Range("V2").Value = 0.00003139085950865095
$ws.Range("W2").Value = 0.00003139085950865095
$ws.Range("M3").Value = '[-0.6935467382578224, 13.82530216424936]'
$ws.Range("N3").Value = 0.07515434006723654
$ws.Range("O3").Value = 0.07515434006723654
$ws.Range("Q3").Value = '[-3.798842768304774, -0.22642109215061534]'
$ws.Range("R3").Value = 0.02808499684366628
$ws.Range("S3").Value = 0.02808499684366628
$ws.Range("U3").Value = '[4.911751158793468, 12.798423435575824]'
$ws.Range("V3").Value = 0.00004420379044334233
$ws.Range("W3").Value = 0.00004420379044334233
$ws.Range("Y3").Value = 0.9246846846847045
$ws.Range("Z3").Value = 15.5141541541545
$ws.Range("M4").Value = '[-1.118666505145077, 14.38740042632763]'
$ws.Range("N4").Value = 0.09166482752672489
$ws.Range("O4").Value = 0.09166482752672489
$ws.Range("Q4").Value = '[-5.5284483333442, 0.6981317007977319]'
$ws.Range("R4").Value = 0.1251866165857611
$ws.Range("S4").Value = 0.1251866165857611
$ws.Range("U4").Value = '[4.786044385430724, 12.688736564651535]'
$ws.Range("V4").Value = 0.00005529420745431501
$ws.Range("W4").Value = 0.00005529420745431501
$ws.Range("Y4").Value = -2.851111111111175
$ws.Range("Z4").Value = 22.57771771771822
$ws.Range("M5").Value = '[-0.5553854873758386, 13.950924173576398]'
$ws.Range("N5").Value = 0.06944659593550062
$ws.Range("O5").Value = 0.06944659593550062
$ws.Range("Q5").Value = '[-4.503263943884468, -1.1069475616252316]'
$ws.Range("R5").Value = 0.001755548272379892
$ws.Range("S5").Value = 0.001755548272379892
$ws.Range("U5").Value = '[4.731233089334403, 12.642138538596331]'
$ws.Range("V5").Value = 0.00006099954326321821
$ws.Range("W5").Value = 0.00006099954326321821
$ws.Range("Y5").Value = 4.520680680680779
$ws.Range("Z5").Value = 18.39095095095136
$ws.Range("M6").Value = '[-0.23228056446629886, 13.565282676843212]'
$ws.Range("N6").Value = 0.05788009967911534
$ws.Range("O6").Value = 0.05788009967911534
$ws.Range("Q6").Value = '[1.6667108172198102, 4.50955341866643]'
$ws.Range("R6").Value = 0.00007106026963987766
$ws.Range("S6").Value = 0.00007106026963987766
$ws.Range("U6").Value = '[4.751498770355218, 12.667930295348272]'
$ws.Range("V6").Value = 0.0000593260762913328
$ws.Range("W6").Value = 0.0000593260762913328
$ws.Range("Y6").Value = 7.243363363363525
$ws.Range("Z6").Value = 18.85329329329371
$ws.Range("M7").Value = '[-1.1255368527570706, 14.289997844899597]'
$ws.Range("N7").Value = 0.09230543459707641
$ws.Range("O7").Value = 0.09230543459707641
$ws.Range("Q7").Value = '[0.8113422468730409, 4.585027116049968]'
$ws.Range("R7").Value = 0.006066359121851272
$ws.Range("S7").Value = 0.006066359121851272
$ws.Range("U7").Value = '[4.851338659845915, 12.755435873279701]'
$ws.Range("V7").Value = 0.00004972439798578421
$ws.Range("W7").Value = 0.00004972439798578421
$ws.Range("Y7").Value = 6.93513513513529
$ws.Range("Z7").Value = 22.34654654654704
$ws.Range("M8").Value = '[-1.8736003466851034, 14.824516478009356]'
$ws.Range("N8").Value = 0.1252654171045577
$ws.Range("O8").Value = 0.1252654171045577
$ws.Range("Q8").Value = '[-0.3522105877898465, 5.497000959434392]'
$ws.Range("R8").Value = 0.08324193924942125
$ws.Range("S8").Value = 0.08324193924942125
$ws.Range("U8").Value = '[5.577872325160936, 14.194785935300096]'
$ws.Range("V8").Value = 0.00003204960463443207
$ws.Range("W8").Value = 0.00003204960463443207
$ws.Range("Y8").Value = 3.210710710710785
$ws.Range("Z8").Value = 27.098398398399
$ws.Range("M9").Value = '[-0.7254561309139564, 14.825614917055061]'
$ws.Range("N9").Value = 0.07446096488810494
$ws.Range("O9").Value = 0.07446096488810494
$ws.Range("Q9").Value = '[0.42139481039142357, 3.3648690083494275]'
$ws.Range("R9").Value = 0.01286133364548681
$ws.Range("S9").Value = 0.01286133364548681
$ws.Range("U9").Value = '[5.225554094695859, 13.41857040591172]'
$ws.Range("V9").Value = 0.00003631420894101289
$ws.Range("W9").Value = 0.00003631420894101289
$ws.Range("Y9").Value = 11.09141141141155
$ws.Range("Z9").Value = 22.27843843843871
$ws.Range("M10").Value = '[-0.5845527772553911, 14.477892204565935]'
$ws.Range("N10").Value = 0.06975142545751667
$ws.Range("O10").Value = 0.06975142545751667
$ws.Range("Q10").Value = '[-0.0503157982556921, 4.415211296937006]'
$ws.Range("R10").Value = 0.05516347885100026
$ws.Range("S10").Value = 0.05516347885100026
$ws.Range("U10").Value = '[5.194951965594104, 13.34226379243255]'
$ws.Range("V10").Value = 0.00003639918568931755
$ws.Range("W10").Value = 0.00003639918568931755
$ws.Range("Y10").Value = 7.099459459459549
$ws.Range("Z10").Value = 24.07123123123153
$ws.Range("M11").Value = '[-0.20668122645954767, 14.063527149803283]'
$ws.Range("N11").Value = 0.05671894895990759
$ws.Range("O11").Value = 0.05671894895990759
$ws.Range("Q11").Value = '[0.7107106503616549, 3.6038690500639667]'
$ws.Range("R11").Value = 0.004346141510672252
$ws.Range("S11").Value = 0.004346141510672252
$ws.Range("U11").Value = '[4.563539513782057, 11.971854441509965]'
$ws.Range("V11").Value = 0.00004830096132901751
$ws.Range("W11").Value = 0.00004830096132901751
$ws.Range("Y11").Value = 10.18306306306319
$ws.Range("Z11").Value = 21.17885885885912
$ws.Range("M12").Value = '[-0.7640581394774575, 14.78930163236312]'
$ws.Range("N12").Value = 0.07600216917268532
$ws.Range("O12").Value = 0.07600216917268532
$ws.Range("Q12").Value = '[0.15723686954903826, 4.107027032620891]'
$ws.Range("R12").Value = 0.03496766915526051
$ws.Range("S12").Value = 0.03496766915526051
$ws.Range("U12").Value = '[5.18856995051746, 13.387705616833218]'
$ws.Range("V12").Value = 0.00003876781926237527
$ws.Range("W12").Value = 0.00003876781926237527
$ws.Range("Y12").Value = 8.27075075075085
$ws.Range("Z12").Value = 23.28240240240269
$ws.Range("M13").Value = '[-0.8269193619803907, 14.90629179335926]'
$ws.Range("N13").Value = 0.07818491448695331
$ws.Range("O13").Value = 0.07818491448695331
$ws.Range("Q13").Value = '[-0.5471843060306547, 5.390079888141046]'
$ws.Range("R13").Value = 0.1073824444003793
$ws.Range("S13").Value = 0.1073824444003793
$ws.Range("U13").Value = '[5.26741313272703, 13.449580644221484]'
$ws.Range("V13").Value = 0.00003358027854183376
$ws.Range("W13").Value = 0.00003358027854183376
$ws.Range("Y13").Value = 3.394354354354398
$ws.Range("Z13").Value = 25.95963963963996
$ws.Range("B14").Value = 0
$ws.Range("M14").Value = '[-0.12842727394414233, 14.000081596088739]'
$ws.Range("N14").Value = 0.05413113984826223
$ws.Range("O14").Value = 0.05413113984826223
$ws.Range("Q14").Value = '[-0.5283158816847697, 5.320895665539468]'
$ws.Range("R14").Value = 0.1058523572811698
$ws.Range("S14").Value = 0.1058523572811698
$ws.Range("U14").Value = '[4.5637243388359625, 12.02632724751295]'
$ws.Range("V14").Value = 0.00005119026329380993
$ws.Range("W14").Value = 0.00005119026329380993
$ws.Range("Y14").Value = 3.657297297297347
$ws.Range("Z14").Value = 25.88792792792825
